$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 4910
$ws.Range("I29").Value = 375
$ws.Range("J29").Value = 7933.3335
$ws.Range("K29").Value = 1125
$ws.Range("L29").Value = 23800.0005
$ws.Range("M29").Value = -844
$ws.Range("N29").Value = -24362.0005

$ws.Range("H76").Value = 4447.3335
$ws.Range("I76").Value = 3417.8
$ws.Range("J76").Value = 9595
$ws.Range("K76").Value = 3417.8
$ws.Range("L76").Value = 9595
$ws.Range("M76").Value = -3102.8

$ws.Range("H79").Value = 4447.3335
$ws.Range("I79").Value = 3417.8
$ws.Range("J79").Value = 9595
$ws.Range("K79").Value = 3417.8
$ws.Range("L79").Value = 9595
$ws.Range("M79").Value = -2325.8

$ws.Range("H132").Value = 818.34283
$ws.Range("I132").Value = 795.24243
$ws.Range("J132").Value = 1199.5
$ws.Range("K132").Value = 2385.72729
$ws.Range("L132").Value = 3598.5
$ws.Range("M132").Value = 144.2727100000002
$ws.Range("N132").Value = -8658.5

$ws.Range("H135").Value = 599.44446
$ws.Range("I135").Value = 526.06665
$ws.Range("J135").Value = 966.3333
$ws.Range("K135").Value = 4734.59985
$ws.Range("L135").Value = 8696.9997
$ws.Range("M135").Value = -2199.59985

$ws.Range("H137").Value = 1557.2069
$ws.Range("I137").Value = 1280.2858
$ws.Range("J137").Value = 2284.125
$ws.Range("K137").Value = 3840.8574
$ws.Range("L137").Value = 6852.375
$ws.Range("M137").Value = -1290.8574
$ws.Range("N137").Value = -11952.375

$ws.Range("H138").Value = 3318.195
$ws.Range("I138").Value = 4734.643
$ws.Range("J138").Value = 2583.7407
$ws.Range("K138").Value = 14203.929
$ws.Range("L138").Value = 7751.222099999999
$ws.Range("M138").Value = -9063.929
$ws.Range("N138").Value = -18031.2221

$ws.Range("H141").Value = 850944.4399999999
$ws.Range("I141").Value = 1168793.9
$ws.Range("J141").Value = 3345.889
$ws.Range("K141").Value = 3506381.7
$ws.Range("L141").Value = 10037.667
$ws.Range("M141").Value = -3501201.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2812.2842
$ws.Range("I32").Value = 2274.2727
$ws.Range("J32").Value = 6578.364
$ws.Range("K32").Value = 2274.2727
$ws.Range("L32").Value = 6578.364
$ws.Range("M32").Value = -1987.2727
$ws.Range("N32").Value = -7152.364

$ws.Range("H37").Value = 13800
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 13800
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 13800
$ws.Range("N37").Value = -14346

$ws.Range("H61").Value = 2344.8696
$ws.Range("I61").Value = 888.13336
$ws.Range("J61").Value = 5076.25
$ws.Range("K61").Value = 888.13336
$ws.Range("L61").Value = 5076.25
$ws.Range("M61").Value = -676.13336

$ws.Range("H74").Value = 1254.6086
$ws.Range("I74").Value = 836.4375
$ws.Range("J74").Value = 2210.4285
$ws.Range("K74").Value = 836.4375
$ws.Range("L74").Value = 2210.4285
$ws.Range("M74").Value = 37.5625

$ws.Range("H77").Value = 1254.6086
$ws.Range("I77").Value = 836.4375
$ws.Range("J77").Value = 2210.4285
$ws.Range("K77").Value = 4182.1875
$ws.Range("L77").Value = 11052.1425
$ws.Range("M77").Value = 185.8125

$ws.Range("H132").Value = 1559.4048
$ws.Range("I132").Value = 1263.8334
$ws.Range("J132").Value = 3332.8333
$ws.Range("K132").Value = 3791.5002
$ws.Range("L132").Value = 9998.499899999999
$ws.Range("M132").Value = -1261.5002

$ws.Range("H136").Value = 2344.8696
$ws.Range("I136").Value = 888.13336
$ws.Range("J136").Value = 5076.25
$ws.Range("K136").Value = 2664.40008
$ws.Range("L136").Value = 15228.75
$ws.Range("M136").Value = -114.4000800000003

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2842.7144
$ws.Range("I105").Value = 2884.4614
$ws.Range("J105").Value = 2300
$ws.Range("K105").Value = 2884.4614
$ws.Range("L105").Value = 2300
$ws.Range("M105").Value = -1137.4614
$ws.Range("N105").Value = -5794

$ws.Range("H132").Value = 124560.375
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 124560.375
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 124560.375
$ws.Range("N132").Value = -134680.375

$ws.Range("H134").Value = 9282.875
$ws.Range("I134").Value = 10052.333
$ws.Range("J134").Value = 3896.6667
$ws.Range("K134").Value = 30156.999
$ws.Range("L134").Value = 11690.0001
$ws.Range("M134").Value = -27621.999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()

$ws.Range("H31").Value = 1384.339
$ws.Range("I31").Value = 697.5
$ws.Range("J31").Value = 1855.3143
$ws.Range("K31").Value = 697.5
$ws.Range("L31").Value = 1855.3143
$ws.Range("M31").Value = -402.5

$ws.Range("H34").Value = 1384.339
$ws.Range("I34").Value = 697.5
$ws.Range("J34").Value = 1855.3143
$ws.Range("K34").Value = 697.5
$ws.Range("L34").Value = 1855.3143
$ws.Range("M34").Value = -495.5

$ws.Range("H58").Value = 2175188.5
$ws.Range("I58").Value = 3106774.8
$ws.Range("J58").Value = 1487.8334
$ws.Range("K58").Value = 3106774.8
$ws.Range("L58").Value = 1487.8334
$ws.Range("M58").Value = -3106571.8
$ws.Range("N58").Value = -1893.8334

$ws.Range("H105").Value = 2000
$ws.Range("I105").Value = 2000
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 2000
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = -253

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()

$ws.Range("H132").Value = 2413.4167
$ws.Range("I132").Value = 1772.5555
$ws.Range("J132").Value = 4336
$ws.Range("K132").Value = 5317.666499999999
$ws.Range("L132").Value = 13008
$ws.Range("M132").Value = -2787.666499999999

$ws.Range("H134").Value = 1693.0465
$ws.Range("I134").Value = 887.4375
$ws.Range("J134").Value = 4036.6365
$ws.Range("K134").Value = 2662.3125
$ws.Range("L134").Value = 12109.9095
$ws.Range("M134").Value = -127.3125
$ws.Range("N134").Value = -17179.9095

$ws.Range("H136").Value = 2175188.5
$ws.Range("I136").Value = 3106774.8
$ws.Range("J136").Value = 1487.8334
$ws.Range("K136").Value = 9320324.399999999
$ws.Range("L136").Value = 4463.5002
$ws.Range("M136").Value = -9317774.399999999
$ws.Range("N136").Value = -9563.5002

$ws.Range("H139").Value = 50000
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 50000
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 50000
$ws.Range("N139").Value = -60280

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 126.57143
$ws.Range("I8").Value = 126.57143
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 379.71429
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -240.71429

$ws.Range("H56").Value = 776189.0600000001
$ws.Range("I56").Value = 776189.0600000001
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 776189.0600000001
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -775659.0600000001

$ws.Range("H68").Value = 3217.2942
$ws.Range("I68").Value = 2002
$ws.Range("J68").Value = 3293.25
$ws.Range("K68").Value = 6006
$ws.Range("L68").Value = 9879.75
$ws.Range("M68").Value = -5195
$ws.Range("N68").Value = -11501.75

$ws.Range("H71").Value = 3217.2942
$ws.Range("I71").Value = 2002
$ws.Range("J71").Value = 3293.25
$ws.Range("K71").Value = 18018
$ws.Range("L71").Value = 29639.25
$ws.Range("M71").Value = -13962
$ws.Range("N71").Value = -37751.25

$ws.Range("H107").Value = 1515.0513
$ws.Range("I107").Value = 1549
$ws.Range("J107").Value = 1510.0588
$ws.Range("K107").Value = 4647
$ws.Range("L107").Value = 4530.1764
$ws.Range("M107").Value = -2727
$ws.Range("N107").Value = -8370.1764

$ws.Range("H109").Value = 4308.3
$ws.Range("I109").Value = 1022.25
$ws.Range("J109").Value = 6499
$ws.Range("K109").Value = 3066.75
$ws.Range("L109").Value = 19497
$ws.Range("M109").Value = -2026.75
$ws.Range("N109").Value = -21577

$ws.Range("H131").Value = 10015059
$ws.Range("I131").Value = 166667100
$ws.Range("J131").Value = 15991.404
$ws.Range("K131").Value = 500001300
$ws.Range("L131").Value = 47974.212
$ws.Range("M131").Value = -499996260
$ws.Range("N131").Value = -58054.212

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 7500
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 7500
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 7500
$ws.Range("N33").Value = -8004

$ws.Range("H113").Value = 1866.6666
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1866.6666
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 1866.6666
$ws.Range("N113").Value = -6206.6666

$ws.Range("H132").Value = 2027120.8
$ws.Range("I132").Value = 2565819.8
$ws.Range("J132").Value = 6999.5
$ws.Range("K132").Value = 7697459.399999999
$ws.Range("L132").Value = 20998.5
$ws.Range("M132").Value = -7694929.399999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 6066.1113
$ws.Range("I32").Value = 5590.625
$ws.Range("J32").Value = 9870
$ws.Range("K32").Value = 5590.625
$ws.Range("L32").Value = 9870
$ws.Range("M32").Value = -5273.625

$ws.Range("H132").Value = 5190.143
$ws.Range("I132").Value = 1450
$ws.Range("J132").Value = 5583.8423
$ws.Range("K132").Value = 4350
$ws.Range("L132").Value = 16751.5269
$ws.Range("M132").Value = -1820
$ws.Range("N132").Value = -21811.5269

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()

$ws.Range("H113").Value = 675.4
$ws.Range("I113").Value = 362.25
$ws.Range("J113").Value = 884.1667
$ws.Range("K113").Value = 1086.75
$ws.Range("L113").Value = 2652.5001
$ws.Range("M113").Value = 1083.25

$ws.Range("H132").Value = 908.5714
$ws.Range("I132").Value = 387.3125
$ws.Range("J132").Value = 2576.6
$ws.Range("K132").Value = 1161.9375
$ws.Range("L132").Value = 7729.799999999999
$ws.Range("M132").Value = 1268.0625
$ws.Range("N132").Value = -12789.8
